$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PathfinderTestAsset")

# Insert a new column before column I (minimum_required_path_nodes),
# shifting existing columns I:P to J:Q.
$ws.Columns("I:I").Insert()

# Populate the header of the newly inserted column.
$ws.Range("I1").Value = "qualifiers"
